$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 300.1111
$ws.Range("I33").Value = 248.14285
$ws.Range("K33").Value = 248.14285
$ws.Range("M33").Value = -19.14285000000001
$ws.Range("H70").Value = 383134.25
$ws.Range("I70").Value = 718454.2
$ws.Range("J70").Value = 3105
$ws.Range("K70").Value = 2155362.6
$ws.Range("L70").Value = 9315
$ws.Range("M70").Value = -2155092.6
$ws.Range("N70").Value = -9855
$ws.Range("H73").Value = 383134.25
$ws.Range("I73").Value = 718454.2
$ws.Range("J73").Value = 3105
$ws.Range("K73").Value = 2155362.6
$ws.Range("L73").Value = 9315
$ws.Range("M73").Value = -2154426.6
$ws.Range("N73").Value = -11187
$ws.Range("H98").Value = 1375.4117
$ws.Range("I98").Value = 1185.6072
$ws.Range("K98").Value = 1185.6072
$ws.Range("M98").Value = 312.3928000000001
$ws.Range("H113").Value = 8435
$ws.Range("I113").Value = 7649.5
$ws.Range("J113").Value = 10006
$ws.Range("K113").Value = 7649.5
$ws.Range("L113").Value = 10006
$ws.Range("M113").Value = -4395.5
$ws.Range("N113").Value = -16514
$ws.Range("H122").Value = 1375.4117
$ws.Range("I122").Value = 1185.6072
$ws.Range("K122").Value = 3556.8216
$ws.Range("M122").Value = -1106.8216
$ws.Range("H127").Value = 3362.3
$ws.Range("I127").Value = 1725
$ws.Range("K127").Value = 5175
$ws.Range("M127").Value = -215
$ws.Range("H132").Value = 4613.6177
$ws.Range("I132").Value = 2290.375
$ws.Range("K132").Value = 6871.125
$ws.Range("M132").Value = -4341.125
$ws.Range("H137").Value = 28504056
$ws.Range("I137").Value = 83336216
$ws.Range("J137").Value = 1087976.4
$ws.Range("K137").Value = 250008648
$ws.Range("L137").Value = 3263929.2
$ws.Range("M137").Value = -250006098
$ws.Range("N137").Value = -3269029.2

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 6955.5557
$ws.Range("I4").Value = 4500
$ws.Range("K4").Value = 4500
$ws.Range("M4").Value = -4384
$ws.Range("H10").Value = 4491
$ws.Range("I10").Value = 5986.5
$ws.Range("K10").Value = 5986.5
$ws.Range("M10").Value = -5816.5
$ws.Range("H26").Value = 12780.7
$ws.Range("I26").Value = 5561.4
$ws.Range("K26").Value = 5561.4
$ws.Range("M26").Value = -5231.4
$ws.Range("H30").Value = 45499
$ws.Range("J30").Value = 89989
$ws.Range("L30").Value = 89989
$ws.Range("N30").Value = -90289
$ws.Range("H32").Value = 1866.6171
$ws.Range("I32").Value = 1675.7046
$ws.Range("K32").Value = 1675.7046
$ws.Range("M32").Value = -1388.7046
$ws.Range("H45").Value = 3766.1667
$ws.Range("I45").Value = 3119.4
$ws.Range("K45").Value = 3119.4
$ws.Range("M45").Value = -2742.4
$ws.Range("H60").Value = 63120.94
$ws.Range("I60").Value = 63120.94
$ws.Range("K60").Value = 63120.94
$ws.Range("M60").Value = -62387.94
$ws.Range("H63").Value = 4997
$ws.Range("I63").Value = 4997
$ws.Range("K63").Value = 4997
$ws.Range("M63").Value = -4311
$ws.Range("H66").Value = 4997
$ws.Range("I66").Value = 4997
$ws.Range("K66").Value = 24985
$ws.Range("M66").Value = -21553
$ws.Range("H88").Value = 3642.1667
$ws.Range("I88").Value = 2406
$ws.Range("K88").Value = 2406
$ws.Range("M88").Value = -2000
$ws.Range("H91").Value = 3642.1667
$ws.Range("I91").Value = 2406
$ws.Range("K91").Value = 2406
$ws.Range("M91").Value = -1002
$ws.Range("H132").Value = 1963968.5
$ws.Range("I132").Value = 3425.2559
$ws.Range("J132").Value = 12501888
$ws.Range("K132").Value = 10275.7677
$ws.Range("L132").Value = 37505664
$ws.Range("M132").Value = -7745.7677
$ws.Range("N132").Value = -37510724

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 2696.5
$ws.Range("I7").Value = 2696.5
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 2696.5
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -2583.5
$ws.Range("N7").ClearContents()
$ws.Range("H86").Value = 25988.758
$ws.Range("J86").Value = 3997.5
$ws.Range("L86").Value = 3997.5
$ws.Range("N86").Value = -6243.5
$ws.Range("H89").Value = 25988.758
$ws.Range("J89").Value = 3997.5
$ws.Range("L89").Value = 19987.5
$ws.Range("N89").Value = -31219.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 4998
$ws.Range("I12").Value = 4998
$ws.Range("K12").Value = 4998
$ws.Range("M12").Value = -4828
$ws.Range("H94").Value = 1314.3889
$ws.Range("I94").Value = 1011.7778
$ws.Range("J94").Value = 1617
$ws.Range("K94").Value = 1011.7778
$ws.Range("L94").Value = 1617
$ws.Range("M94").Value = -560.7778
$ws.Range("N94").Value = -2519
$ws.Range("H111").Value = 99999
$ws.Range("J111").Value = 99999
$ws.Range("L111").Value = 99999
$ws.Range("N111").Value = -108179
$ws.Range("H132").Value = 2872.5
$ws.Range("I132").Value = 2764.2
$ws.Range("J132").Value = 3414
$ws.Range("K132").Value = 8292.599999999999
$ws.Range("L132").Value = 10242
$ws.Range("M132").Value = -5762.599999999999
$ws.Range("N132").Value = -15302
$ws.Range("H134").Value = 3321
$ws.Range("I134").Value = 3384.1538
$ws.Range("K134").Value = 10152.4614
$ws.Range("M134").Value = -7617.4614

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5255.875
$ws.Range("I3").Value = 1244.8572
$ws.Range("K3").Value = 3734.5716
$ws.Range("M3").Value = -3622.5716
$ws.Range("H131").Value = 4486.7666
$ws.Range("I131").Value = 3594.875
$ws.Range("J131").Value = 4811.091
$ws.Range("K131").Value = 10784.625
$ws.Range("L131").Value = 14433.273
$ws.Range("M131").Value = -5744.625
$ws.Range("N131").Value = -24513.273

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2421.1667
$ws.Range("I102").Value = 2421.1667
$ws.Range("K102").Value = 2421.1667
$ws.Range("M102").Value = -799.1667000000002
$ws.Range("H132").Value = 4549742
$ws.Range("I132").Value = 3139.25
$ws.Range("K132").Value = 9417.75
$ws.Range("M132").Value = -6887.75

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2808.25
$ws.Range("I40").Value = 2280.8572
$ws.Range("J40").Value = 6500
$ws.Range("K40").Value = 2280.8572
$ws.Range("L40").Value = 6500
$ws.Range("M40").Value = -2144.8572
$ws.Range("N40").Value = -6772
$ws.Range("H68").Value = 9724141
$ws.Range("I68").Value = 10804157
$ws.Range("J68").Value = 4003
$ws.Range("K68").Value = 10804157
$ws.Range("L68").Value = 4003
$ws.Range("M68").Value = -10803408
$ws.Range("N68").Value = -5501
$ws.Range("H71").Value = 9724141
$ws.Range("I71").Value = 10804157
$ws.Range("J71").Value = 4003
$ws.Range("K71").Value = 54020785
$ws.Range("L71").Value = 20015
$ws.Range("M71").Value = -54017041
$ws.Range("N71").Value = -27503
$ws.Range("H93").Value = 3273413
$ws.Range("I93").Value = 2776.2222
$ws.Range("J93").Value = 6952879
$ws.Range("K93").Value = 2776.2222
$ws.Range("L93").Value = 6952879
$ws.Range("M93").Value = -1528.2222
$ws.Range("N93").Value = -6955375
$ws.Range("H132").Value = 6480.6665
$ws.Range("I132").Value = 3333
$ws.Range("J132").Value = 7110.2
$ws.Range("K132").Value = 9999
$ws.Range("L132").Value = 21330.6
$ws.Range("M132").Value = -7469
$ws.Range("N132").Value = -26390.6

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 10279.857
$ws.Range("I62").Value = 5492.25
$ws.Range("K62").Value = 5492.25
$ws.Range("M62").Value = -4868.25
$ws.Range("H65").Value = 10279.857
$ws.Range("I65").Value = 5492.25
$ws.Range("K65").Value = 27461.25
$ws.Range("M65").Value = -24341.25
$ws.Range("H132").Value = 834390.5600000001
$ws.Range("I132").Value = 1068.7
$ws.Range("K132").Value = 3206.1
$ws.Range("M132").Value = -676.1000000000004
$ws.Range("H135").Value = 110000
$ws.Range("J135").Value = 110000
$ws.Range("L135").Value = 110000
$ws.Range("N135").Value = -120140
